$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap "Periodo Mora" (E) and "Valor Mora" (F) between rows 16 and 17
$e16 = $ws.Range("E16").Value2
$e17 = $ws.Range("E17").Value2
$f16 = $ws.Range("F16").Value2
$f17 = $ws.Range("F17").Value2

$ws.Range("E16").Value2 = $e17
$ws.Range("E17").Value2 = $e16
$ws.Range("F16").Value2 = $f17
$ws.Range("F17").Value2 = $f16
